$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("G2").Value = 181.4944075
$ws.Range("H2").Value = 362.988815
$ws.Range("I2").Value = 0.2239486468210351
$ws.Range("J2").Value = 0.1654349085470023
$ws.Range("M2").Value = 0.4333629999999999
$ws.Range("N2").Value = 0.8667259999999999
$ws.Range("Q2").Value = 78.6529609174225
$ws.Range("R2").Value = 314.61184366969
$ws.Range("S2").Value = 0.2239486468210351
$ws.Range("T2").Value = 0.1654349085470023

# Row 3
$ws.Range("I3").Value = 0.07700606288633029
$ws.Range("J3").Value = 0.08532865336765341
$ws.Range("M3").Value = 0.4333629999999999
$ws.Range("N3").Value = 0.8667259999999999
$ws.Range("Q3").Value = 27.04528444614033
$ws.Range("R3").Value = 162.271706676842
$ws.Range("S3").Value = 0.07700606288633029
$ws.Range("T3").Value = 0.08532865336765341

# Row 4
$ws.Range("G4").Value = 171.9980316666667
$ws.Range("H4").Value = 515.994095
$ws.Range("I4").Value = 0.2122309275432167
$ws.Range("J4").Value = 0.235168226649403
$ws.Range("M4").Value = 0.4333629999999999
$ws.Range("N4").Value = 0.8667259999999999
$ws.Range("Q4").Value = 74.53758299716165
$ws.Range("R4").Value = 447.22549798297
$ws.Range("S4").Value = 0.2122309275432167
$ws.Range("T4").Value = 0.235168226649403

# Row 5
$ws.Range("G5").Value = 55.64279550000001
$ws.Range("H5").Value = 111.285591
$ws.Range("I5").Value = 0.06865847234198982
$ws.Range("J5").Value = 0.05071925307032974
$ws.Range("M5").Value = 0.4333629999999999
$ws.Range("N5").Value = 0.8667259999999999
$ws.Range("Q5").Value = 24.1135287862665
$ws.Range("R5").Value = 96.454115145066
$ws.Range("S5").Value = 0.06865847234198982
$ws.Range("T5").Value = 0.05071925307032974

# Row 6
$ws.Range("G6").Value = 203.386317
$ws.Range("H6").Value = 610.158951
$ws.Range("I6").Value = 0.250961399315095
$ws.Range("J6").Value = 0.2780845747487284
$ws.Range("M6").Value = 0.4333629999999999
$ws.Range("N6").Value = 0.8667259999999999
$ws.Range("Q6").Value = 88.14010449407098
$ws.Range("R6").Value = 528.840626964426
$ws.Range("S6").Value = 0.250961399315095
$ws.Range("T6").Value = 0.2780845747487284

# Row 7
$ws.Range("G7").Value = 135.4992116666667
$ws.Range("H7").Value = 406.497635
$ws.Range("I7").Value = 0.167194491092333
$ws.Range("J7").Value = 0.1852643836168829
$ws.Range("M7").Value = 0.4333629999999999
$ws.Range("N7").Value = 0.8667259999999999
$ws.Range("Q7").Value = 58.72034486550166
$ws.Range("R7").Value = 352.32206919301
$ws.Range("S7").Value = 0.167194491092333
$ws.Range("T7").Value = 0.1852643836168829
